# This workbook is a single-sheet SDG 10.3.1 indicator table with three
# language columns (A = Kyrgyz, B = Russian, C = English). Two section
# header rows are re-worded:
#   Row 19 - "Age (in years)" breakdown header
#   Row 29 - "Education" breakdown header
#
# The wording changes from a bare noun ("Age (in years)" / "Education") to
# a "By ..." / "По ..." / "... боюнча" phrasing in all three languages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A text (Kyrgyz) for the two section headers
$ws.Range("A19").Value = "Жаш курагы боюнча (жылдарда)"   # was "Жаш курагы (жылдарда)"
$ws.Range("A29").Value = "Билими боюнча"                   # was "Билими"

# Column B text (Russian) for the two section headers
$ws.Range("B19").Value = "По возрасту (в годах)"           # was "Возраст (в годах)"
$ws.Range("B29").Value = "По образованию"                  # was "Образование"

# Column C text (English) for the two section headers
$ws.Range("C19").Value = "By age (in years) "              # was "Age (in years) "
$ws.Range("C29").Value = "By education"                    # was "Education"
